# ADD results from server
# Updates row 2 (data row) values on sheets "2025", "2030", and "2035"
# to reflect newly computed results from the server.

$wb = $excel.ActiveWorkbook

# --- Sheet "2025" ---
$ws = $wb.Worksheets.Item([string]"2025")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.0008630959698206382
$ws.Range("E2").Value = 0.3707083019056102
$ws.Range("I2").Value = 0.6522940196752842
$ws.Range("L2").Value = 0.3150122101148058
$ws.Range("M2").Value = 0.08569991666666667
$ws.Range("N2").Value = 12.81572300722258
$ws.Range("O2").Value = 3.087566746787787

# --- Sheet "2030" ---
$ws = $wb.Worksheets.Item([string]"2030")
$ws.Range("A2").Value = 0.006109625212652015
$ws.Range("B2").Value = 0.04996758146600301
$ws.Range("E2").Value = 0.2221320845210674
$ws.Range("I2").Value = 0.5435695833333334
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0.04727991666666674
$ws.Range("N2").Value = 5.250308734235212
$ws.Range("O2").Value = 2.141760504764548

# --- Sheet "2035" ---
$ws = $wb.Worksheets.Item([string]"2035")
$ws.Range("A2").Value = 0.2127705
$ws.Range("B2").Value = 0.04380969999999995
$ws.Range("E2").Value = 0.08036855622576544
$ws.Range("I2").Value = 0.4209396056630867
$ws.Range("M2").Value = 0.04866624788015655
$ws.Range("N2").Value = 4.304825345210443
$ws.Range("O2").Value = 6.77691883767349
